$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add X3 and Y3 (new "Down" sentiment classification columns)
$ws.Range("X3").Value = -1.3299870000000169
$ws.Range("Y3").Value = "Down"

# Row 4: new data row appended to the sheet
$ws.Range("A4").Value = 42633.888333333336
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 16924
$ws.Range("F4").Value = 2955
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 41
$ws.Range("I4").Value = 82
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 30419
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = 260
$ws.Range("N4").Value = 44
$ws.Range("O4").Value = 9
$ws.Range("P4").Value = "Noun"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 1.76
$ws.Range("S4").Value = 0.1055
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = -6.67
$ws.Range("U4").Value = 5.83
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
